# updated outputs, changed axis to 0.01
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D5 was the text "not clear" (a shared string); change it to the numeric value 0.8
$ws.Range("D5").Value = 0.8

# Update the active selection in the sheet view from D11 to H11
$ws.Range("H11").Select()
